$d = $word.ActiveDocument

# 1. Append " (Questions 1-5)" to the bold note paragraph, then add a
#    trailing space (kept bold) after it.
$notePara = $d.Paragraphs(4)
$noteRange = $notePara.Range
$noteRange.MoveEnd(1, -1) | Out-Null
$noteRange.Text = "Please note that the steps show rounded numbers, but that the final answers to the problems are calculated without rounding. (Questions 1-5)"

$notePara2 = $d.Paragraphs(4)
$noteRange2 = $notePara2.Range
$insertPoint = $d.Range($noteRange2.End - 1, $noteRange2.End - 1)
$insertPoint.InsertAfter(" ")
$insertPoint.Font.Bold = $true

# 2. The table had a duplicated "Mode" row/value baked in by a double
#    file error. Remove the standalone "1 / Mode / The most frequently
#    occurring value" row entirely, and strip the trailing
#    " Mode -NN.NNNNNN" text that had been appended to each of the
#    Company A-E solution cells.
$t = $d.Tables(1)

$modeRow = $t.Rows(4)
if ($modeRow.Cells(2).Range.Text -like "Mode*") {
    $modeRow.Delete()
}

$replacements = @(
    "Mean: 21.276 Median: 13.433",
    "Mean: 33.482 Median: 20.838",
    "Mean: 41.122 Median: 25.558",
    "Mean: 0.706 Median: 1.892",
    "Mean: -1.084 Median: -3.796"
)
for ($i = 0; $i -lt 5; $i++) {
    $row = $t.Rows(4 + $i)
    $cell = $row.Cells(3)
    $cellRange = $cell.Range
    $cellRange.MoveEnd(1, -1) | Out-Null
    $cellRange.Text = $replacements[$i]
}
